# Atualização de bases das ligas, do dia: 24-02-2024 às 12:40
# Re-sync a handful of match rows (ids re-shuffled on re-scrape) in the
# "Denmark Superligaen" sheet. Column A (sequence id) stays put; every
# other column (B..AC) for the affected rows takes the values belonging
# to its corresponding match record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 26, 27, 28: 3-way cyclic rotation of the match data
#   new row26 = old row27, new row27 = old row28, new row28 = old row26
# ---------------------------------------------------------------------
$row26 = $ws.Range("B26:AC26").Value2
$row27 = $ws.Range("B27:AC27").Value2
$row28 = $ws.Range("B28:AC28").Value2

$ws.Range("B26:AC26").Value2 = $row27
$ws.Range("B27:AC27").Value2 = $row28
$ws.Range("B28:AC28").Value2 = $row26

# ---------------------------------------------------------------------
# Rows 30 <-> 31: simple swap
# ---------------------------------------------------------------------
$row30 = $ws.Range("B30:AC30").Value2
$row31 = $ws.Range("B31:AC31").Value2

$ws.Range("B30:AC30").Value2 = $row31
$ws.Range("B31:AC31").Value2 = $row30

# ---------------------------------------------------------------------
# Rows 33 <-> 34: simple swap
# ---------------------------------------------------------------------
$row33 = $ws.Range("B33:AC33").Value2
$row34 = $ws.Range("B34:AC34").Value2

$ws.Range("B33:AC33").Value2 = $row34
$ws.Range("B34:AC34").Value2 = $row33

# ---------------------------------------------------------------------
# Rows 190 <-> 191: simple swap
# ---------------------------------------------------------------------
$row190 = $ws.Range("B190:AC190").Value2
$row191 = $ws.Range("B191:AC191").Value2

$ws.Range("B190:AC190").Value2 = $row191
$ws.Range("B191:AC191").Value2 = $row190

# ---------------------------------------------------------------------
# Rows 203-205: odds re-calculation (individual cell corrections)
# ---------------------------------------------------------------------
$ws.Range("N203").Value2 = 2.9
$ws.Range("P203").Value2 = 2.375
$ws.Range("R203").Value2 = 1.81
$ws.Range("S203").Value2 = 2.09
$ws.Range("T203").Value2 = 2.25
$ws.Range("U203").Value2 = 1.8
$ws.Range("V203").Value2 = 2.05

$ws.Range("R204").Value2 = 1.89
$ws.Range("S204").Value2 = 2.01
$ws.Range("U204").Value2 = 1.9
$ws.Range("V204").Value2 = 1.95

$ws.Range("N205").Value2 = 2.9
$ws.Range("O205").Value2 = 3.4
$ws.Range("P205").Value2 = 2.4

Write-Output "Applied league base update for 2024-02-24 12:40"
